$wb = $excel.ActiveWorkbook

# RB sheet gets the new Week 16 log entry for L.Bell
$ws = $wb.Worksheets.Item("RB")

$ws.Range("A6").Value = "L.Bell"
$ws.Range("B6:J6").Value = 0

# Make RB the active sheet/tab with the selection left at H7 (ready for next week's entry)
$ws.Activate()
$ws.Range("H7").Select()
